$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Rename "ENGLISH" -> "English"
# ---------------------------------------------------------------------------
$wsEnglish = $wb.Worksheets.Item("ENGLISH")
$wsEnglish.Name = "English"

# ---------------------------------------------------------------------------
# 2. Touch every populated cell on the pre-existing sheets so they pick up
#    the (new) explicit "Normal" cell style - mirrors the style-index bump
#    (s="0" -> s="1") seen across tags / notes / English in the target file.
#    Iterating cell-by-cell (instead of whole rows/columns) avoids
#    materialising phantom empty cells in sparsely populated rows.
# ---------------------------------------------------------------------------
$wsTags = $wb.Worksheets.Item("tags")
for ($r = 1; $r -le 8; $r++) {
    $cell = $wsTags.Cells.Item($r, 1)
    if ($cell.Text -ne "") {
        $cell.Style = "Normal"
    }
}

$wsNotes = $wb.Worksheets.Item("notes")
for ($r = 1; $r -le 49; $r++) {
    for ($c = 1; $c -le 3; $c++) {
        $cell = $wsNotes.Cells.Item($r, $c)
        if ($cell.Text -ne "") {
            $cell.Style = "Normal"
        }
    }
}

for ($r = 1; $r -le 85; $r++) {
    for ($c = 1; $c -le 7; $c++) {
        $cell = $wsEnglish.Cells.Item($r, $c)
        if ($cell.Text -ne "") {
            $cell.Style = "Normal"
        }
    }
}

# ---------------------------------------------------------------------------
# 3. Add the new "settings" sheet at the end of the workbook and populate it.
# ---------------------------------------------------------------------------
$wsSettings = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$wsSettings.Name = "settings"

$wsSettings.Range("A1").Value = "default language name"
$wsSettings.Range("B1").Value = "entries per vocabulary training session"
$wsSettings.Range("A2").Value = "English"
$wsSettings.Range("B2").Value = 6

$wsSettings.Columns.Item(1).ColumnWidth = 18.76
$wsSettings.Columns.Item(2).ColumnWidth = 36.54

Write-Output "done"
